# Apply updates to the "electric-targets" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("electric-targets")

# Update the electric_target (column C) values for the "Articulated trucks"
# rows (rows 20-35, sales_year 2025-2040).
$ws.Range("C20").Value = 2
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 7
$ws.Range("C23").Value = 11
$ws.Range("C24").Value = 15
$ws.Range("C25").Value = 20
$ws.Range("C26").Value = 25
$ws.Range("C27").Value = 30
$ws.Range("C28").Value = 35
$ws.Range("C29").Value = 40
$ws.Range("C30").Value = 45
$ws.Range("C31").Value = 50
$ws.Range("C32").Value = 55
$ws.Range("C33").Value = 60
$ws.Range("C34").Value = 65
$ws.Range("C35").Value = 70

# Update the share_applied (column D) values for every data row (2-35) from
# 0.8 to 1.
$ws.Range("D2:D35").Value = 1

# Update the selected cell shown when the workbook is reopened.
$ws.Range("G33").Select()
